# Update gh-pages output data (scraped counters) across sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

$oldImg = "//i1.hdslb.com/bfs/openplatform/202406/suoZa5Ha1717727447336.jpeg"
$newImg = "//i1.hdslb.com/bfs/openplatform/202406/TnP82LF01719390282936.jpeg"

# --- 展览 (Exhibition) sheet ---
$wsExhibit.Range("F3").Value  = 760
$wsExhibit.Range("F4").Value  = 1512
$wsExhibit.Range("F6").Value  = 95
$wsExhibit.Range("F7").Value  = 154
$wsExhibit.Range("I7").Value  = $newImg
$wsExhibit.Range("F8").Value  = 6260
$wsExhibit.Range("F10").Value = 405
$wsExhibit.Range("F12").Value = 5248
$wsExhibit.Range("F13").Value = 29
$wsExhibit.Range("F17").Value = 59
$wsExhibit.Range("F18").Value = 365
$wsExhibit.Range("F19").Value = 70
$wsExhibit.Range("F21").Value = 302
$wsExhibit.Range("F23").Value = 3743
$wsExhibit.Range("F24").Value = 164

# --- 演出 (Show) sheet ---
$wsShow.Range("F2").Value = 85

# --- 全部类型 (All types) sheet ---
$wsAll.Range("F2").Value  = 85
$wsAll.Range("F4").Value  = 760
$wsAll.Range("F5").Value  = 1512
$wsAll.Range("F7").Value  = 95
$wsAll.Range("F8").Value  = 154
$wsAll.Range("I8").Value  = $newImg
$wsAll.Range("F9").Value  = 6260
$wsAll.Range("F11").Value = 405
$wsAll.Range("F13").Value = 5248
$wsAll.Range("F14").Value = 29
$wsAll.Range("F18").Value = 59
$wsAll.Range("F19").Value = 365
$wsAll.Range("F20").Value = 70
$wsAll.Range("F22").Value = 302
$wsAll.Range("F24").Value = 3743
$wsAll.Range("F26").Value = 164
